$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - Home Odds / Away Odds
$ws.Range("E11").Value = 1.2
$ws.Range("F11").Value = 4.1

# Row 54 - Home Odds / Away Odds
$ws.Range("E54").Value = 1.92
$ws.Range("F54").Value = 1.8

# Row 56 - Home Odds / Away Odds
$ws.Range("E56").Value = 2.22
$ws.Range("F56").Value = 1.6
